$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/number-like readings coming from meteo.cat's daily scrape.
# Assigning these as plain strings keeps them as literal text (General
# format, same style as before) for every cell except H2: "95%" would
# otherwise be auto-recognised as a percentage number by Excel's normal
# typed-entry parsing, so that one is entered with a leading apostrophe
# (the standard Excel "force text" quote-prefix) to keep it literal text.
$ws.Range("E2").Value = "2026-02-20 09:15:36"
$ws.Range("H2").Value = "'95%"
$ws.Range("I2").Value = "0.1 mm"
$ws.Range("J2").Value = "1020.5 hPa"
$ws.Range("K2").Value = "1.6 MJ/m2"
$ws.Range("M2").Value = "10.3 °C 8:59 TU"
$ws.Range("O2").Value = "2.3 °C"
